$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Job to Run")

# Update job name value in A2 (was "TJGIp11", now "TJGIp11_pos")
$ws.Range("A2").Value = "TJGIp11_pos"

# Update the active cell/selection on the sheet
$ws.Range("D11").Select()

$wb.Save()
